$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 132
$ws.Range("F5").Value = 364
$ws.Range("F6").Value = 760
$ws.Range("F7").Value = 204
$ws.Range("F8").Value = 1074
$ws.Range("F9").Value = 276
$ws.Range("F12").Value = 625
$ws.Range("F13").Value = 170
$ws.Range("F14").Value = 494
$ws.Range("F17").Value = 156
$ws.Range("F18").Value = 825
$ws.Range("F19").Value = 2592
$ws.Range("F20").Value = 521
$ws.Range("F23").Value = 307
$ws.Range("F24").Value = 206
$ws.Range("F26").Value = 151
$ws.Range("F28").Value = 962
$ws.Range("F29").Value = 4
$ws.Range("F31").Value = 227
$ws.Range("F32").Value = 1022
$ws.Range("F34").Value = 44

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 1044
$ws.Range("F5").Value = 1044
$ws.Range("F14").Value = 585
$ws.Range("F15").Value = 91
$ws.Range("F16").Value = 10
$ws.Range("F19").Value = 35
$ws.Range("F20").Value = 610
$ws.Range("F22").Value = 34
$ws.Range("F24").Value = 291
$ws.Range("F25").Value = 262
$ws.Range("F26").Value = 3697
$ws.Range("F28").Value = 4
$ws.Range("F30").Value = 192
$ws.Range("F31").Value = 23
$ws.Range("F33").Value = 106
$ws.Range("F35").Value = 5

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1768
$ws.Range("F5").Value = 2401
$ws.Range("F6").Value = 986
$ws.Range("F9").Value = 1243
$ws.Range("F10").Value = 327
$ws.Range("F11").Value = 89

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1768
$ws.Range("F4").Value = 2401
$ws.Range("F7").Value = 986
$ws.Range("F8").Value = 1243
$ws.Range("F9").Value = 327
$ws.Range("F10").Value = 89
$ws.Range("F11").Value = 132
$ws.Range("F12").Value = 364
$ws.Range("F13").Value = 760
$ws.Range("F14").Value = 204
$ws.Range("F16").Value = 1074
$ws.Range("F17").Value = 276
$ws.Range("F19").Value = 625
$ws.Range("F20").Value = 1044
$ws.Range("F21").Value = 494
$ws.Range("F23").Value = 156
$ws.Range("F24").Value = 825
$ws.Range("F25").Value = 2592
$ws.Range("F26").Value = 521
$ws.Range("F28").Value = 307
$ws.Range("F30").Value = 206
$ws.Range("F31").Value = 151
$ws.Range("F34").Value = 962
$ws.Range("F35").Value = 585
$ws.Range("F36").Value = 585
$ws.Range("F37").Value = 91
$ws.Range("F39").Value = 227
$ws.Range("F40").Value = 10
$ws.Range("F41").Value = 35
$ws.Range("F43").Value = 34
$ws.Range("F44").Value = 291
$ws.Range("F45").Value = 291
$ws.Range("F46").Value = 262
$ws.Range("F47").Value = 1022
$ws.Range("F49").Value = 192
$ws.Range("F50").Value = 44
